$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 285, shifting the existing data (and all rows below)
# down by two rows.
$ws.Rows.Item(285).Resize(2).Insert()

# Row 285: new record
$ws.Cells.Item(285, 1).Value = 7
$ws.Cells.Item(285, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(285, 3).Value = "Ñuble"
$ws.Cells.Item(285, 4).Value = 44809
$ws.Cells.Item(285, 5).Value = 16
$ws.Cells.Item(285, 6).Value = 100114001
$ws.Cells.Item(285, 7).Value = "Papa"
$ws.Cells.Item(285, 8).Value = "Patagonia"
$ws.Cells.Item(285, 9).Value = "1a (guarda)"
$ws.Cells.Item(285, 10).Value = 120
$ws.Cells.Item(285, 11).Value = 7000
$ws.Cells.Item(285, 12).Value = 7500
$ws.Cells.Item(285, 13).Value = 7250
$ws.Cells.Item(285, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(285, 15).Value = "Región de Ñuble"
$ws.Cells.Item(285, 16).Value = 290
$ws.Cells.Item(285, 17).Value = 25
$ws.Cells.Item(285, 18).Value = "Hortaliza"

# Row 286: new record
$ws.Cells.Item(286, 1).Value = 7
$ws.Cells.Item(286, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(286, 3).Value = "Ñuble"
$ws.Cells.Item(286, 4).Value = 44809
$ws.Cells.Item(286, 5).Value = 16
$ws.Cells.Item(286, 6).Value = 100114001
$ws.Cells.Item(286, 7).Value = "Papa"
$ws.Cells.Item(286, 8).Value = "Rosara"
$ws.Cells.Item(286, 9).Value = "1a (guarda)"
$ws.Cells.Item(286, 10).Value = 120
$ws.Cells.Item(286, 11).Value = 7500
$ws.Cells.Item(286, 12).Value = 8000
$ws.Cells.Item(286, 13).Value = 7750
$ws.Cells.Item(286, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(286, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(286, 16).Value = 310
$ws.Cells.Item(286, 17).Value = 25
$ws.Cells.Item(286, 18).Value = "Hortaliza"

# Apply the same date format style as other date cells in column D
$ws.Range("D285:D286").NumberFormat = $ws.Range("D287").NumberFormat
